$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" header) values per regenerated save_data
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("G6").Value = 3
